$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column widths for columns A and B (15.42578125 -> 14.42578125 chars).
# The COM ColumnWidth setter quantizes to whole pixels, so we feed it the
# character width (13.666666666666666) whose nearest-pixel rounding lands
# on the target stored width (14.42578125 rounds to the 14.5 bucket).
$ws.Columns("A:B").ColumnWidth = 13.666666666666666

# Update cell values (row 4 is unchanged)
$ws.Range("A1").Value = -0.032510287630593084
$ws.Range("B1").Value = 0.032510286740648264

$ws.Range("A2").Value = 0.010791408273078798
$ws.Range("B2").Value = -0.010791409188199368

$ws.Range("A3").Value = 0.0042568731580106111
$ws.Range("B3").Value = -0.0042568740807759986
